# Atualização de bases das ligas, do dia: 17-06-2024 às 21:10
#
# The source feed re-ordered a handful of fixture rows (column A keeps the
# running rank, but every other field - match id, teams, odds, results...
# - belongs to the other row of the pair). Swap columns B:AD between each
# pair of rows so the data lines back up with the correct rank/index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (1-based worksheet rows) whose B:AD contents must be exchanged.
$pairs = @(
    @(69, 70),
    @(117, 118),
    @(120, 121),
    @(161, 162),
    @(184, 185)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1" + ":AD$r1")
    $range2 = $ws.Range("B$r2" + ":AD$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
